$d = $word.ActiveDocument

# 1. Update the first paragraph's text.
$d.Content.Find.Execute("**************", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Changed demo1 and demo3 document.", 2)

# 2. Insert a new paragraph right after the (now updated) first paragraph.
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()

# 3. Grab the newly created (now second) paragraph - it already inherited the
#    same tab stop formatting from paragraph 1 - and set its text.
$p2 = $d.Paragraphs(2)
$p2.Range.Text = "Demo 2 document remains same"
